# The rows of weekly Fruta/Hortaliza price data (rows 2-9, excluding row 7)
# were reshuffled: each row's Fecha/Volumen/Precio/Unidad de comercializacion/
# Precio $/Kg values were rotated among the rows. Row 7 is untouched.
# Apply the new values directly, cell by cell, matching the target XML.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (previously row 6's data)
$ws.Range("D2").Value = 44309
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 7000
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 7000
$ws.Range("Q2").Value = "$/caja 14 kilos empedrada"
$ws.Range("S2").Value = 500

# Row 3 (previously row 5's data)
$ws.Range("D3").Value = 44400
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 10000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 10000
$ws.Range("Q3").Value = "$/caja 14 kilos"
$ws.Range("S3").Value = 714

# Row 4 (previously row 2's data)
$ws.Range("D4").Value = 44176
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 7000
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 7000
$ws.Range("Q4").Value = "$/caja 14 kilos empedrada"
$ws.Range("S4").Value = 500

# Row 5 (previously row 9's data)
$ws.Range("D5").Value = 44397
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 11000
$ws.Range("O5").Value = 11000
$ws.Range("P5").Value = 11000
$ws.Range("Q5").Value = "$/caja 14 kilos"
$ws.Range("S5").Value = 786

# Row 6 (previously row 8's data)
$ws.Range("D6").Value = 44491
$ws.Range("M6").Value = 180
$ws.Range("N6").Value = 9000
$ws.Range("O6").Value = 9000
$ws.Range("P6").Value = 9000
$ws.Range("Q6").Value = "$/caja 14 kilos empedrada"
$ws.Range("S6").Value = 643

# Row 8 (previously row 3's data)
$ws.Range("D8").Value = 44208
$ws.Range("M8").Value = 210
$ws.Range("N8").Value = 10000
$ws.Range("O8").Value = 10000
$ws.Range("P8").Value = 10000
$ws.Range("Q8").Value = "$/caja 14 kilos empedrada"
$ws.Range("S8").Value = 714

# Row 9 (previously row 4's data)
$ws.Range("D9").Value = 44351
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 10000
$ws.Range("Q9").Value = "$/caja 14 kilos empedrada"
$ws.Range("S9").Value = 714
